# Daily refresh of the cryptos price table (rows 2-51 of the active sheet).
# Price/volume figures are updated, and rows 7/8 (USDC <-> LidoStakedEther)
# swap places in the ranking.
#
# Price cells in column D are plain text (they use "." as a thousands
# separator, e.g. "60.887.64"), so numeric-looking values are written with a
# leading apostrophe to force Excel to keep them as text instead of
# reinterpreting them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.887.64'
$ws.Range('E2').Value = '  -0.80%  '
$ws.Range('D3').Value = '3.387.52'
$ws.Range('E3').Value = '  -1.36%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('D5').Value = '''571.38'
$ws.Range('E5').Value = '  -0.75%  '
$ws.Range('D6').Value = '''141.72'
$ws.Range('E6').Value = '  -2.32%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D8').Value = '3.388.02'
$ws.Range('E8').Value = '  -1.38%  '
$ws.Range('E9').Value = '  -0.21%  '
$ws.Range('E10').Value = '  -1.79%  '
$ws.Range('E11').Value = '  -2.07%  '
$ws.Range('D12').Value = '''0.394'
$ws.Range('E12').Value = '  +2.13%  '
$ws.Range('D13').Value = '3.965.29'
$ws.Range('E13').Value = '  -1.37%  '
$ws.Range('E14').Value = '  +2.16%  '
$ws.Range('D15').Value = '''28.19'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('E16').Value = '  -1.28%  '
$ws.Range('D17').Value = '3.387.88'
$ws.Range('E17').Value = '  -1.28%  '
$ws.Range('D18').Value = '60.957.61'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').Value = '''6.15'
$ws.Range('E19').Value = '  -2.18%  '
$ws.Range('D20').Value = '''13.84'
$ws.Range('E20').Value = '  -2.86%  '
$ws.Range('E21').Value = '  -4.53%  '
$ws.Range('D22').Value = '''383.82'
$ws.Range('E22').Value = '  -2.92%  '
$ws.Range('D23').Value = '''0.557'
$ws.Range('E23').Value = '  -1.53%  '
$ws.Range('D24').Value = '''74.42'
$ws.Range('E24').Value = '  +0.94%  '
$ws.Range('E25').Value = '  +0.61%  '
$ws.Range('D26').Value = '''0.0000117'
$ws.Range('E26').Value = '  -4.63%  '
$ws.Range('D27').Value = '3.527.11'
$ws.Range('E27').Value = '  -1.33%  '
$ws.Range('E28').Value = '  -1.23%  '
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('E30').Value = '  -2.82%  '
$ws.Range('E31').Value = '  -3.21%  '
$ws.Range('E32').Value = '  -1.87%  '
$ws.Range('E33').Value = '  -2.46%  '
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D35').Value = '''23.48'
$ws.Range('E35').Value = '  -1.88%  '
$ws.Range('D36').Value = '''6.97'
$ws.Range('E36').Value = '  -0.47%  '
$ws.Range('D37').Value = '''167.54'
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('D38').Value = '3.416.35'
$ws.Range('E38').Value = '  -1.36%  '
$ws.Range('D39').Value = '''4.97'
$ws.Range('E39').Value = '  -2.64%  '
$ws.Range('E40').Value = '  -4.37%  '
$ws.Range('E41').Value = '  -1.43%  '
$ws.Range('D42').Value = '''27.54'
$ws.Range('E42').Value = '  +2.14%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').Value = '''0.779'
$ws.Range('E44').Value = '  -2.49%  '
$ws.Range('D45').Value = '''42.14'
$ws.Range('E45').Value = '  -0.30%  '
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('E47').Value = '  -3.53%  '
$ws.Range('E48').Value = '  -1.02%  '
$ws.Range('D49').Value = '2.477.47'
$ws.Range('E49').Value = '  -4.38%  '
$ws.Range('D50').Value = '''6.81'
$ws.Range('E50').Value = '  -1.55%  '
$ws.Range('D51').Value = '''23.00'
$ws.Range('E51').Value = '  -1.23%  '
